$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the two new cells on the existing row 5 (PriceChange, UpDown)
$ws.Range("X5").Value = -0.59999799999999937
$ws.Range("Y5").Value = "Down"

# Append a brand new row 6 with a full set of scan results
$ws.Range("A6").Value = 42647.883229166669
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = "Buy"
$ws.Range("D6").Value = 12
$ws.Range("E6").Value = 28624
$ws.Range("F6").Value = 2977
$ws.Range("G6").Value = 50
$ws.Range("H6").Value = 43
$ws.Range("I6").Value = 78
$ws.Range("J6").Value = 21
$ws.Range("K6").Value = 12617
$ws.Range("L6").Value = 387
$ws.Range("M6").Value = 335
$ws.Range("N6").Value = 82
$ws.Range("O6").Value = 22
$ws.Range("P6").Value = "Bag"
$ws.Range("Q6").Value = 53.235658945584888
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = -0.086199999999999999
$ws.Range("T6").Value = -0.0166
$ws.Range("U6").Value = 6.69
$ws.Range("V6").Value = 1.88
$ws.Range("W6").Value = 0

# Match the date/percentage number formats used by the rows above
# (copy format only so we reuse the existing style entries instead of
# minting brand-new ones).
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)

$ws.Range("S5").Copy()
$ws.Range("S6").PasteSpecial(-4122)

$ws.Range("T5").Copy()
$ws.Range("T6").PasteSpecial(-4122)

$excel.CutCopyMode = 0
